$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain looking numbers (e.g. "1.00", "567.70") as
# literal text in the source workbook. A direct .Value assignment lets Excel
# auto-detect these as numbers and silently drop the formatting (trailing
# zeros, etc.), so each Price cell is briefly switched to the Text number
# format while its new value is written, then restored to the default/Normal
# style so no stray cell formatting is introduced.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.764.45'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.661.54'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.613'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.659.74'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.146'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.103.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.683.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.660.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '342.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.38'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0800'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '158.08'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.23'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.09'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.910'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.901'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '303.13'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0988'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.605'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0544'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.33'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E2').Value = '  +3.74%  '
$ws.Range('E3').Value = '  +1.76%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  +6.21%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +5.78%  '
$ws.Range('E9').Value = '  +1.49%  '
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('E11').Value = '  +4.59%  '
$ws.Range('E12').Value = '  +6.77%  '
$ws.Range('E13').Value = '  +3.42%  '
$ws.Range('E14').Value = '  +0.85%  '
$ws.Range('E15').Value = '  +3.68%  '
$ws.Range('E16').Value = '  +5.53%  '
$ws.Range('E17').Value = '  +4.23%  '
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('E19').Value = '  +2.90%  '
$ws.Range('E20').Value = '  +2.38%  '
$ws.Range('E21').Value = '  +2.90%  '
$ws.Range('E22').Value = '  +2.48%  '
$ws.Range('E23').Value = '  +0.80%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('E26').Value = '  +5.20%  '
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +4.14%  '
$ws.Range('E30').Value = '  +9.02%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  +4.45%  '
$ws.Range('E33').Value = '  +3.63%  '
$ws.Range('E34').Value = '  +2.83%  '
$ws.Range('E35').Value = '  +1.50%  '
$ws.Range('E36').Value = '  +4.75%  '
$ws.Range('E37').Value = '  +8.95%  '
$ws.Range('E38').Value = '  +5.18%  '
$ws.Range('E39').Value = '  +10.48%  '
$ws.Range('E40').Value = '  +1.04%  '
$ws.Range('E41').Value = '  +6.51%  '
$ws.Range('E42').Value = '  +5.91%  '
$ws.Range('E43').Value = '  +1.73%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('E45').Value = '  +4.80%  '
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('E47').Value = '  +3.16%  '
$ws.Range('E48').Value = '  +13.81%  '
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('E51').Value = '  +5.38%  '

Write-Output "Applied 86 cell updates"
